# sigGOterms.xlsx edit: bioperl code to extract bed file elements
#
# 1) Nudge the saved window X position (cosmetic; best-effort).
# 2) On sheet "sigGOterms.txt" (2nd sheet), populate column D ("GO_terms"
#    detail text) for a handful of rows that previously had no entry.
# 3) Add a new column F ("sig") with a flag formula
#    =IF(D{row}<>"",1,0) for every data row, filled down in the same
#    batches the author used (so Excel's shared-formula grouping matches).
# 4) Update the active selection to H173 (no frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 1) window position -----------------------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 1020
} catch {
}

# --- 2) fill in newly-annotated GO term descriptions in column D -------
$dValues = @(
    @{Row=6;   Val="transcription factor TFTC complex"},
    @{Row=15;  Val="tRNA aminoacylation for protein translation"},
    @{Row=25;  Val="ER to Golgi vesicle"},
    @{Row=32;  Val="negative regulation of phosphorylation"},
    @{Row=39;  Val="armadillo repeat domain binding"},
    @{Row=42;  Val="histone deacetylation complex"},
    @{Row=46;  Val="jak-stat cascade"},
    @{Row=106; Val="response to bacterium"},
    @{Row=114; Val="vesicle transport along microtuble"},
    @{Row=123; Val="protein serine/theonine phosphotase activity"},
    @{Row=163; Val="desmosome"},
    @{Row=166; Val="small-subunit processome"},
    @{Row=169; Val="alpha tubulin binding"}
)

foreach ($d in $dValues) {
    $ws.Cells.Item($d.Row, 4).Value = $d.Val
}

# --- 3) new column F: header + fill-down formulas -----------------------
$ws.Range("F1").Value = "sig"

$ws.Range("F2:F5").Formula = "=IF(D2<>"""",1,0)"
$ws.Range("F6").Formula = "=IF(D6<>"""",1,0)"
$ws.Range("F7:F70").Formula = "=IF(D7<>"""",1,0)"
$ws.Range("F71:F134").Formula = "=IF(D71<>"""",1,0)"
$ws.Range("F135:F181").Formula = "=IF(D135<>"""",1,0)"

# --- 4) selection update -------------------------------------------------
$ws.Range("H173").Select()
